$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.087.79'
$ws.Range('E2').Value = '  +3.17%  '
$ws.Range('D3').Value = '2.298.06'
$ws.Range('E3').Value = '  +2.04%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.27'
$ws.Range('E5').Value = '  +2.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.03'
$ws.Range('E6').Value = '  +7.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.537'
$ws.Range('E7').Value = '  +2.67%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('E9').Value = '  +8.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.92'
$ws.Range('E10').Value = '  +4.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0824'
$ws.Range('E11').Value = '  +5.06%  '
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('E13').Value = '  +7.98%  '
$ws.Range('D14').Value = '2.654.88'
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.01'
$ws.Range('E15').Value = '  +5.44%  '
$ws.Range('D16').Value = '2.297.77'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.807'
$ws.Range('E17').Value = '  +3.10%  '
$ws.Range('D18').Value = '43.015.98'
$ws.Range('E18').Value = '  +3.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.54'
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('D20').Value = '0.0₃0925'
$ws.Range('E20').Value = '  +3.47%  '
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.52'
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.39'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.62'
$ws.Range('E24').Value = '  +2.96%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  +5.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.47'
$ws.Range('E27').Value = '  +4.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.52'
$ws.Range('E28').Value = '  +7.33%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  +4.36%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.64'
$ws.Range('E30').Value = '  +2.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.60'
$ws.Range('E31').Value = '  +5.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.33'
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.12'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.77'
$ws.Range('E35').Value = '  +5.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0741'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('E37').Value = '  +3.60%  '
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.116'
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.22'
$ws.Range('E41').Value = '  +6.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.29'
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0289'
$ws.Range('E43').Value = '  +3.21%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.973.69'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.01'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.04'
$ws.Range('E46').Value = '  +4.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.86'
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.88'
$ws.Range('E48').Value = '  +6.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.95'
$ws.Range('E49').Value = '  +17.23%  '
$ws.Range('D50').Value = '2.523.11'
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('E51').Value = '  +3.14%  '
